$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'254.68"
$ws.Range("D2").Style = $ws.Range("B2").Style
$ws.Range("E2").Value = "'3.57%"
$ws.Range("E2").Style = $ws.Range("B2").Style

$ws.Range("D3").Value = "'27.96"
$ws.Range("D3").Style = $ws.Range("B3").Style
$ws.Range("E3").Value = "'-7.18%"
$ws.Range("E3").Style = $ws.Range("B3").Style

$ws.Range("D4").Value = "'5.212"
$ws.Range("D4").Style = $ws.Range("B4").Style
$ws.Range("E4").Value = "'1.20%"
$ws.Range("E4").Style = $ws.Range("B4").Style

$ws.Range("D5").Value = "'0.05869"
$ws.Range("D5").Style = $ws.Range("B5").Style
$ws.Range("E5").Value = "'1.89%"
$ws.Range("E5").Style = $ws.Range("B5").Style

$ws.Range("D6").Value = "'6.715"
$ws.Range("D6").Style = $ws.Range("B6").Style
$ws.Range("E6").Value = "'0.70%"
$ws.Range("E6").Style = $ws.Range("B6").Style

$ws.Range("B7").Value = 'GateToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D7").Value = "'3.219"
$ws.Range("D7").Style = $ws.Range("B7").Style
$ws.Range("E7").Value = "'-1.60%"
$ws.Range("E7").Style = $ws.Range("B7").Style

$ws.Range("B8").Value = 'MXToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D8").Value = "'0.8637"
$ws.Range("D8").Style = $ws.Range("B8").Style
$ws.Range("E8").Value = "'1.69%"
$ws.Range("E8").Style = $ws.Range("B8").Style

$ws.Range("B9").Value = 'FTXToken'
$ws.Range("C9").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D9").Value = "'0.9706"
$ws.Range("D9").Style = $ws.Range("B9").Style
$ws.Range("E9").Value = "'13.09%"
$ws.Range("E9").Style = $ws.Range("B9").Style

$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D10").Value = "'0.1411"
$ws.Range("D10").Style = $ws.Range("B10").Style
$ws.Range("E10").Value = "'1.92%"
$ws.Range("E10").Style = $ws.Range("B10").Style

$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D11").Value = "'0.07172"
$ws.Range("D11").Style = $ws.Range("B11").Style
$ws.Range("E11").Value = "'1.26%"
$ws.Range("E11").Style = $ws.Range("B11").Style

$ws.Range("B12").Value = 'BitrueCoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D12").Value = "'0.03177"
$ws.Range("D12").Style = $ws.Range("B12").Style
$ws.Range("E12").Value = "'-1.83%"
$ws.Range("E12").Style = $ws.Range("B12").Style

$ws.Range("B13").Value = 'BitMartToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D13").Value = "'0.09218"
$ws.Range("D13").Style = $ws.Range("B13").Style
$ws.Range("E13").Value = "'-1.51%"
$ws.Range("E13").Style = $ws.Range("B13").Style

$ws.Range("B14").Value = 'BitForexToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D14").Value = "'0.001546"
$ws.Range("D14").Style = $ws.Range("B14").Style
$ws.Range("E14").Value = "'1.18%"
$ws.Range("E14").Style = $ws.Range("B14").Style

$ws.Range("B15").Value = 'One'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D15").Value = "'0.0006049"
$ws.Range("D15").Style = $ws.Range("B15").Style
$ws.Range("E15").Value = "'1.34%"
$ws.Range("E15").Style = $ws.Range("B15").Style

$ws.Range("B16").Value = 'TigerCash'
$ws.Range("C16").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D16").Value = "'0.005795"
$ws.Range("D16").Style = $ws.Range("B16").Style
$ws.Range("E16").Value = "'-1.85%"
$ws.Range("E16").Style = $ws.Range("B16").Style

$ws.Range("B17").Value = 'LEO'
$ws.Range("C17").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D17").Value = "'3.501"
$ws.Range("D17").Style = $ws.Range("B17").Style
$ws.Range("E17").Value = "'-1.30%"
$ws.Range("E17").Style = $ws.Range("B17").Style

$ws.Range("D18").Value = "'2.222"
$ws.Range("D18").Style = $ws.Range("B18").Style
$ws.Range("E18").Value = "'0.27%"
$ws.Range("E18").Style = $ws.Range("B18").Style

$ws.Range("E19").Value = "'1.85%"
$ws.Range("E19").Style = $ws.Range("B19").Style

$ws.Range("D20").Value = "'0.03478"
$ws.Range("D20").Style = $ws.Range("B20").Style
$ws.Range("E20").Value = "'1.69%"
$ws.Range("E20").Style = $ws.Range("B20").Style

$ws.Range("E21").Value = "'-1.73%"
$ws.Range("E21").Style = $ws.Range("B21").Style

$ws.Range("D22").Value = "'3.561"
$ws.Range("D22").Style = $ws.Range("B22").Style
$ws.Range("E22").Value = "'2.00%"
$ws.Range("E22").Style = $ws.Range("B22").Style

$ws.Range("E23").Value = "'1.34%"
$ws.Range("E23").Style = $ws.Range("B23").Style

$ws.Range("E24").Value = "'-2.07%"
$ws.Range("E24").Style = $ws.Range("B24").Style

$ws.Range("D25").Value = "'0.001226"
$ws.Range("D25").Style = $ws.Range("B25").Style
$ws.Range("E25").Value = "'0.02%"
$ws.Range("E25").Style = $ws.Range("B25").Style

$ws.Range("D26").Value = "'0.004798"
$ws.Range("D26").Style = $ws.Range("B26").Style
$ws.Range("E26").Value = "'15.42%"
$ws.Range("E26").Style = $ws.Range("B26").Style

$ws.Range("D27").Value = "'0.0001200"
$ws.Range("D27").Style = $ws.Range("B27").Style
$ws.Range("E27").Value = "'0.02%"
$ws.Range("E27").Style = $ws.Range("B27").Style

$ws.Range("E28").Value = "'1.17%"
$ws.Range("E28").Style = $ws.Range("B28").Style

$ws.Range("E40").Value = "'1.53%"
$ws.Range("E40").Style = $ws.Range("B40").Style

$ws.Range("B41").Value = 'KickToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range("D41").Value = "'0.005732"
$ws.Range("D41").Style = $ws.Range("B41").Style
$ws.Range("E41").Value = "'1.26%"
$ws.Range("E41").Style = $ws.Range("B41").Style

$ws.Range("B42").Value = 'BKEXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range("D42").Value = "'0.1101"
$ws.Range("D42").Style = $ws.Range("B42").Style
$ws.Range("E42").Value = "'2.92%"
$ws.Range("E42").Style = $ws.Range("B42").Style

$ws.Range("D43").Value = "'0.002338"
$ws.Range("D43").Style = $ws.Range("B43").Style
$ws.Range("E43").Value = "'6.28%"
$ws.Range("E43").Style = $ws.Range("B43").Style

$ws.Range("D44").Value = "'0.009462"
$ws.Range("D44").Style = $ws.Range("B44").Style
$ws.Range("E44").Value = "'0.20%"
$ws.Range("E44").Style = $ws.Range("B44").Style

$ws.Range("D45").Value = "'0.00005230"
$ws.Range("D45").Style = $ws.Range("B45").Style
$ws.Range("E45").Value = "'-4.50%"
$ws.Range("E45").Style = $ws.Range("B45").Style

$ws.Range("E46").Value = "'-0.03%"
$ws.Range("E46").Style = $ws.Range("B46").Style

$ws.Range("D47").Value = "'0.09998"
$ws.Range("D47").Style = $ws.Range("B47").Style
$ws.Range("E47").Value = "'40.83%"
$ws.Range("E47").Style = $ws.Range("B47").Style

$ws.Range("D48").Value = "'0.002134"
$ws.Range("D48").Style = $ws.Range("B48").Style
$ws.Range("E48").Value = "'-13.47%"
$ws.Range("E48").Style = $ws.Range("B48").Style

$ws.Range("E49").Value = "'-0.03%"
$ws.Range("E49").Style = $ws.Range("B49").Style

$ws.Range("E50").Value = "'-0.03%"
$ws.Range("E50").Style = $ws.Range("B50").Style
